$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Continue the "Cedula" numeric sequence starting at row 6 (117100595) through row 37 (117100626)
$startRow = 6
$startVal = 117100595
$endRow = 37

for ($r = $startRow; $r -le $endRow; $r++) {
    $val = $startVal + ($r - $startRow)
    $ws.Cells.Item($r, 1).Value = $val
}

# Row 10 previously held an (empty) underline style in A10; now it just holds the value with default style.
$ws.Cells.Item(10, 1).Font.Underline = $false

# Update the selected cell to B3, matching the saved selection in the workbook view.
$ws.Range("B3").Select()
